$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.998.13"
$ws.Range("E2").Value = "  +6.72%  "

$ws.Range("D3").Value = "1.739.83"
$ws.Range("E3").Value = "  +5.06%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.19%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "229.17"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.50%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5471"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.40%  "

$ws.Range("E7").Value = "  -0.24%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2787"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.58%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06730"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.87%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.83"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.58%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07800"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.97%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.710"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.55%  "

$ws.Range("D16").Value = "0.0₅8435"
$ws.Range("E16").Value = "  +2.24%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "69.67"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +6.53%  "

$ws.Range("D18").Value = "27.972.58"
$ws.Range("E18").Value = "  +6.62%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "226.69"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +17.79%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.851"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.35%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.002"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.15%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.98"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.39%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.250"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.28%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.003"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.22%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.49"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.90%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1254"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.13%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.479"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.81%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.16"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +7.70%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.655"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +9.25%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05687"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.78%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.318"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.16%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.706"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.53%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.551"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.74%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.669"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.42%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9858"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.51%  "

$ws.Range("E36").Value = "  +2.09%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.450"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.44%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5975"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.81%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01678"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.91%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.016"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.22%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.002"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.19%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "102.27"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.34%  "

$ws.Range("D45").Value = "1.883.99"
$ws.Range("E45").Value = "  +4.94%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.324"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.95%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.010"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.30%  "

$ws.Range("E50").Value = "  +1.79%  "

$ws.Range("E51").Value = "  -0.67%  "

$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "1.978.43"
$ws.Range("E13").Value = "  +5.02%  "

$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6025"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +7.06%  "

$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "1.677.90"
$ws.Range("E15").Value = "  +4.11%  "

$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "1.049.72"
$ws.Range("E41").Value = "  +3.64%  "

$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8466"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.44%  "

$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "0.0₈117"
$ws.Range("E46").Value = "  +10.59%  "

$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "60.27"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.21%  "
